$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "http://localhost:3000"
$ws.Range("E2").Value = "Galleta Casino Clásica"

$ws.Range("B3").Value = "http://localhost:3000"
$ws.Range("E3").Value = "Galleta Casino Coco"
$ws.Range("G3").Value = "Operación cancelada"
